$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.04585075083389
$ws.Range("C2").Value = 7.360301762512989
$ws.Range("D2").Value = 6.375954411399011
$ws.Range("E2").Value = 11.22302240966287
$ws.Range("F2").Value = 33.16469232162416
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 26.03093852400601
$ws.Range("K2").Value = 11.69856910908128
$ws.Range("M2").Value = 15.07257873823612
$ws.Range("N2").Value = 20.84667234438747
$ws.Range("B3").Value = 10.80082906769157
$ws.Range("C3").Value = 7.133769545368759
$ws.Range("D3").Value = 6.380980532588322
$ws.Range("E3").Value = 11.00936252128352
$ws.Range("F3").Value = 33.03285283925478
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 26.03275416418871
$ws.Range("K3").Value = 11.52957008164125
$ws.Range("M3").Value = 14.91521557112782
$ws.Range("N3").Value = 20.89536727138235
$ws.Range("B4").Value = 10.65085100784465
$ws.Range("C4").Value = 6.993376444673416
$ws.Range("D4").Value = 6.384089450744809
$ws.Range("E4").Value = 10.8798142531005
$ws.Range("F4").Value = 32.95988283099749
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 26.03889924467063
$ws.Range("K4").Value = 11.42787674318554
$ws.Range("M4").Value = 14.82201611513336
$ws.Range("N4").Value = 20.92715163685552
$ws.Range("B5").Value = 10.58994961172887
$ws.Range("C5").Value = 6.935939230037841
$ws.Range("D5").Value = 6.385362302922623
$ws.Range("E5").Value = 10.82751168497956
$ws.Range("F5").Value = 32.93217201055987
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 26.04266711946319
$ws.Range("K5").Value = 11.38701150192754
$ws.Range("M5").Value = 14.78493821984134
$ws.Range("N5").Value = 20.94057831919088
$ws.Range("B6").Value = 10.57985294228719
$ws.Range("C6").Value = 6.926391260935559
$ws.Range("D6").Value = 6.385574023689522
$ws.Range("E6").Value = 10.81885869001654
$ws.Range("F6").Value = 32.92769344303189
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 26.04336905068478
$ws.Range("K6").Value = 11.38026219553631
$ws.Range("M6").Value = 14.77883708267218
$ws.Range("N6").Value = 20.94283646129544
$ws.Range("B7").Value = 10.65002866439173
$ws.Range("C7").Value = 6.992602602544845
$ws.Range("D7").Value = 6.384106592574565
$ws.Range("E7").Value = 10.87910679910252
$ws.Range("F7").Value = 32.95950089161897
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 26.0389449449432
$ws.Range("K7").Value = 11.42732321847624
$ws.Range("M7").Value = 14.82151236614231
$ws.Range("N7").Value = 20.92733079318943
$ws.Range("B8").Value = 10.96132898094371
$ws.Range("C8").Value = 7.282522458804448
$ws.Range("D8").Value = 6.377682791312759
$ws.Range("E8").Value = 11.14905838544905
$ws.Range("F8").Value = 33.11758968645424
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 26.03051996356074
$ws.Range("K8").Value = 11.63990096958488
$ws.Range("M8").Value = 15.01763489671001
$ws.Range("N8").Value = 20.86307112705746
$ws.Range("B9").Value = 11.57128493093081
$ws.Range("C9").Value = 7.836600268301091
$ws.Range("D9").Value = 6.365257980602368
$ws.Range("E9").Value = 11.68807213035458
$ws.Range("F9").Value = 33.48998794826809
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 26.05394463608835
$ws.Range("K9").Value = 12.07078432531332
$ws.Range("M9").Value = 15.42748970908411
$ws.Range("N9").Value = 20.75201045912927
$ws.Range("B10").Value = 12.01373217978522
$ws.Range("C10").Value = 8.229830199128955
$ws.Range("D10").Value = 6.356221315688481
$ws.Range("E10").Value = 12.08552788236733
$ws.Range("F10").Value = 33.80020555910967
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 26.09552222169432
$ws.Range("K10").Value = 12.39263563988687
$ws.Range("M10").Value = 15.74137279413268
$ws.Range("N10").Value = 20.67951803343161
$ws.Range("B11").Value = 12.21272592008276
$ws.Range("C11").Value = 8.40479135985964
$ws.Range("D11").Value = 6.352127407666846
$ws.Range("E11").Value = 12.26575487718646
$ws.Range("F11").Value = 33.94892998160799
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 26.11971907741258
$ws.Range("K11").Value = 12.53951656960185
$ws.Range("M11").Value = 15.88635823635127
$ws.Range("N11").Value = 20.64851411875027
$ws.Range("B12").Value = 12.28767135265407
$ws.Range("C12").Value = 8.470413656892479
$ws.Range("D12").Value = 6.350579371116578
$ws.Range("E12").Value = 12.33384801025654
$ws.Range("F12").Value = 34.00630835720625
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 26.12963952796233
$ws.Range("K12").Value = 12.59514833144282
$ws.Range("M12").Value = 15.9415274915184
$ws.Range("N12").Value = 20.6370574388409
$ws.Range("B13").Value = 12.27154995313977
$ws.Range("C13").Value = 8.456309822894889
$ws.Range("D13").Value = 6.350912672058935
$ws.Range("E13").Value = 12.31919093300028
$ws.Range("F13").Value = 33.99390432991131
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 26.1274693254747
$ws.Range("K13").Value = 12.5831674203294
$ws.Range("M13").Value = 15.92963474674578
$ws.Range("N13").Value = 20.63951221494228
$ws.Range("B14").Value = 12.21890037476049
$ws.Range("C14").Value = 8.410203169928252
$ws.Range("D14").Value = 6.352000006038806
$ws.Range("E14").Value = 12.27136044431339
$ws.Range("F14").Value = 33.9536294992811
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 26.12052006899685
$ws.Range("K14").Value = 12.54409352177279
$ws.Range("M14").Value = 15.89089198920226
$ws.Range("N14").Value = 20.64756588190275
$ws.Range("B15").Value = 12.18659537043947
$ws.Range("C15").Value = 8.381877352659052
$ws.Range("D15").Value = 6.352666315120602
$ws.Range("E15").Value = 12.24204064014908
$ws.Range("F15").Value = 33.92909696109533
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 26.11636204184828
$ws.Range("K15").Value = 12.52015943706397
$ws.Range("M15").Value = 15.86719415273363
$ws.Range("N15").Value = 20.65253595117035
$ws.Range("B16").Value = 12.00067453443097
$ws.Range("C16").Value = 8.218311214408192
$ws.Range("D16").Value = 6.35648918591551
$ws.Range("E16").Value = 12.07373178104736
$ws.Range("F16").Value = 33.79063624652447
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 26.09404704614768
$ws.Range("K16").Value = 12.38304121657602
$ws.Range("M16").Value = 15.73193764151433
$ws.Range("N16").Value = 20.6815839252436
$ws.Range("B17").Value = 11.88597494150143
$ws.Range("C17").Value = 8.116913904412778
$ws.Range("D17").Value = 6.358838586993137
$ws.Range("E17").Value = 11.9702796125653
$ws.Range("F17").Value = 33.7076194115527
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 26.08170935970577
$ws.Range("K17").Value = 12.29900324166704
$ws.Range("M17").Value = 15.64948984815154
$ws.Range("N17").Value = 20.69990936259458
$ws.Range("B18").Value = 11.81979292431559
$ws.Range("C18").Value = 8.058227384880203
$ws.Range("D18").Value = 6.360191504110305
$ws.Range("E18").Value = 11.91072638856319
$ws.Range("F18").Value = 33.66058829691567
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 26.07511047416149
$ws.Range("K18").Value = 12.25071426368483
$ws.Range("M18").Value = 15.60227797015992
$ws.Range("N18").Value = 20.71063539514311
$ws.Range("B19").Value = 11.79735136555714
$ws.Range("C19").Value = 8.03829654975589
$ws.Range("D19").Value = 6.36064986002314
$ws.Range("E19").Value = 11.89055630057129
$ws.Range("F19").Value = 33.6447886814061
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 26.07296169058808
$ws.Range("K19").Value = 12.2343743347574
$ws.Range("M19").Value = 15.58633040763114
$ws.Range("N19").Value = 20.71429894206963
$ws.Range("B20").Value = 11.89820723574966
$ws.Range("C20").Value = 8.1277461593296
$ws.Range("D20").Value = 6.358588324649907
$ws.Range("E20").Value = 11.98129800695614
$ws.Range("F20").Value = 33.71638262682124
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 26.08297125865548
$ws.Range("K20").Value = 12.30794473033931
$ws.Range("M20").Value = 15.65824519176658
$ws.Range("N20").Value = 20.69793936605076
$ws.Range("B21").Value = 12.23437656082089
$ws.Range("C21").Value = 8.423763462058115
$ws.Range("D21").Value = 6.351680570512888
$ws.Range("E21").Value = 12.28541417822954
$ws.Range("F21").Value = 33.96543071460216
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 26.1225406919414
$ws.Range("K21").Value = 12.55557061151389
$ws.Range("M21").Value = 15.90226484147329
$ws.Range("N21").Value = 20.64519262089702
$ws.Range("B22").Value = 12.45166364124577
$ws.Range("C22").Value = 8.613519207751352
$ws.Range("D22").Value = 6.347178906249054
$ws.Range("E22").Value = 12.48323786679784
$ws.Range("F22").Value = 34.13435745252261
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 26.15281621964313
$ws.Range("K22").Value = 12.71744844515387
$ws.Range("M22").Value = 16.06327574614066
$ws.Range("N22").Value = 20.61237400298425
$ws.Range("B23").Value = 12.33594033661308
$ws.Range("C23").Value = 8.512602981347941
$ws.Range("D23").Value = 6.349580408270223
$ws.Range("E23").Value = 12.37776391638431
$ws.Range("F23").Value = 34.04364620360418
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 26.13625451100947
$ws.Range("K23").Value = 12.63106571894821
$ws.Range("M23").Value = 15.9772175853352
$ws.Range("N23").Value = 20.62973850494506
$ws.Range("B24").Value = 11.89267775581735
$ws.Range("C24").Value = 8.12285011221071
$ws.Range("D24").Value = 6.358701461303211
$ws.Range("E24").Value = 11.97631682584299
$ws.Range("F24").Value = 33.71241860535039
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 26.08239921441788
$ws.Range("K24").Value = 12.30390219931331
$ws.Range("M24").Value = 15.65428631009758
$ws.Range("N24").Value = 20.69882940802636
$ws.Range("B25").Value = 11.4069011802228
$ws.Range("C25").Value = 7.688803054119189
$ws.Range("D25").Value = 6.368602145547875
$ws.Range("E25").Value = 11.54168738711931
$ws.Range("F25").Value = 33.38270164401249
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 26.04332712664095
$ws.Range("K25").Value = 11.9530562154802
$ws.Range("M25").Value = 15.31417686407633
$ws.Range("N25").Value = 20.78045583283116
